$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values like "398.32" would otherwise be auto-converted to numbers by Excel)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '56.308.65'
$ws.Range('E2').Value = '  +10.00%  '
$ws.Range('D3').Value = '3.229.98'
$ws.Range('E3').Value = '  +5.07%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '398.32'
$ws.Range('E5').Value = '  +3.00%  '
$ws.Range('D6').Value = '111.49'
$ws.Range('E6').Value = '  +8.60%  '
$ws.Range('D7').Value = '0.556'
$ws.Range('E7').Value = '  +3.59%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.621'
$ws.Range('E9').Value = '  +6.88%  '
$ws.Range('D10').Value = '39.29'
$ws.Range('E10').Value = '  +6.86%  '
$ws.Range('D11').Value = '0.0930'
$ws.Range('E11').Value = '  +9.31%  '
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('D13').Value = '3.734.96'
$ws.Range('E13').Value = '  +4.35%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '19.17'
$ws.Range('E14').Value = '  +4.35%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '8.10'
$ws.Range('E15').Value = '  +5.22%  '
$ws.Range('D16').Value = '3.237.19'
$ws.Range('E16').Value = '  +4.50%  '
$ws.Range('E17').Value = '  +5.47%  '
$ws.Range('D18').Value = '11.11'
$ws.Range('E18').Value = '  +4.81%  '
$ws.Range('D19').Value = '56.162.98'
$ws.Range('E19').Value = '  +9.41%  '
$ws.Range('E20').Value = '  +4.37%  '
$ws.Range('E21').Value = '  +8.15%  '
$ws.Range('D22').Value = '13.07'
$ws.Range('E22').Value = '  +5.44%  '
$ws.Range('D23').Value = '298.50'
$ws.Range('E23').Value = '  +12.77%  '
$ws.Range('D24').Value = '75.98'
$ws.Range('E24').Value = '  +8.72%  '
$ws.Range('E25').Value = '  +2.55%  '
$ws.Range('D26').Value = '8.19'
$ws.Range('E26').Value = '  +3.63%  '
$ws.Range('D27').Value = '28.16'
$ws.Range('E27').Value = '  +4.27%  '
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('E29').Value = '  +4.43%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +4.75%  '
$ws.Range('D32').Value = '11.18'
$ws.Range('E32').Value = '  +6.67%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = '37.02'
$ws.Range('E33').Value = '  +3.87%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').Value = '0.0491'
$ws.Range('E34').Value = '  +3.89%  '
$ws.Range('E35').Value = '  +3.07%  '
$ws.Range('D36').Value = '51.45'
$ws.Range('E36').Value = '  +3.44%  '
$ws.Range('D37').Value = '3.54'
$ws.Range('E37').Value = '  +5.28%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '3.11'
$ws.Range('E38').Value = '  +25.94%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').Value = '136.84'
$ws.Range('E40').Value = '  +6.00%  '
$ws.Range('D41').Value = '17.53'
$ws.Range('E41').Value = '  +6.28%  '
$ws.Range('D42').Value = '1.93'
$ws.Range('E42').Value = '  +4.94%  '
$ws.Range('D43').Value = '4.00'
$ws.Range('E43').Value = '  +5.01%  '
$ws.Range('D44').Value = '0.120'
$ws.Range('E44').Value = '  +3.84%  '
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('D46').Value = '22.26'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('D47').Value = '2.21'
$ws.Range('E47').Value = '  +56.40%  '
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D49').Value = '2.132.46'
$ws.Range('E49').Value = '  +3.35%  '
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('D51').Value = '0.0361'
$ws.Range('E51').Value = '  +11.45%  '
